$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 19-22 (values shift to reflect newer weekly records) ---

# Row 19
$ws.Range("D19").Value = 44417
$ws.Range("M19").Value = 56
$ws.Range("N19").Value = 16000
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 16000
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("T19").Value = 10

# Row 20
$ws.Range("D20").Value = 44417
$ws.Range("L20").Value = "Segunda"
$ws.Range("M20").Value = 60
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 14000
$ws.Range("P20").Value = 14000
$ws.Range("S20").Value = 1400

# Row 21
$ws.Range("D21").Value = 44420
$ws.Range("M21").Value = 54
$ws.Range("N21").Value = 18000
$ws.Range("O21").Value = 18000
$ws.Range("P21").Value = 18000
$ws.Range("S21").Value = 1800

# Row 22
$ws.Range("D22").Value = 44420
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 15000
$ws.Range("O22").Value = 15000
$ws.Range("P22").Value = 15000
$ws.Range("S22").Value = 1500

# Row 23 - now carries what used to be the (older) first record, dated earliest
$ws.Range("D23").Value = 44319
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 60
$ws.Range("N23").Value = 24000
$ws.Range("O23").Value = 24000
$ws.Range("P23").Value = 24000
$ws.Range("Q23").Value = "$/caja 15 kilos granel"
$ws.Range("R23").Value = "Provincia del Elquí"
$ws.Range("S23").Value = 1600
$ws.Range("T23").Value = 15

# --- Append new rows 24-27 with the remaining historical weekly records ---

# Row 24
$ws.Range("A24").Value = 3
$ws.Range("B24").Value = "Femacal de La Calera"
$ws.Range("C24").Value = "Coquimbo"
$ws.Range("D24").Value = 44370
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100108
$ws.Range("H24").Value = "Tropicales y subtropicales"
$ws.Range("I24").Value = 100108004
$ws.Range("J24").Value = "Papaya"
$ws.Range("K24").Value = "Cultivar IV Región"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 17000
$ws.Range("O24").Value = 17000
$ws.Range("P24").Value = 17000
$ws.Range("Q24").Value = "$/bandeja 10 kilos"
$ws.Range("R24").Value = "Provincia del Elquí"
$ws.Range("S24").Value = 1700
$ws.Range("T24").Value = 10

# Row 25
$ws.Range("A25").Value = 3
$ws.Range("B25").Value = "Femacal de La Calera"
$ws.Range("C25").Value = "Coquimbo"
$ws.Range("D25").Value = 44382
$ws.Range("E25").Value = 5
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100108
$ws.Range("H25").Value = "Tropicales y subtropicales"
$ws.Range("I25").Value = 100108004
$ws.Range("J25").Value = "Papaya"
$ws.Range("K25").Value = "Cultivar IV Región"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 58
$ws.Range("N25").Value = 17000
$ws.Range("O25").Value = 17000
$ws.Range("P25").Value = 17000
$ws.Range("Q25").Value = "$/bandeja 10 kilos"
$ws.Range("R25").Value = "Provincia del Elquí"
$ws.Range("S25").Value = 1700
$ws.Range("T25").Value = 10

# Row 26
$ws.Range("A26").Value = 3
$ws.Range("B26").Value = "Femacal de La Calera"
$ws.Range("C26").Value = "Coquimbo"
$ws.Range("D26").Value = 44398
$ws.Range("E26").Value = 5
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100108
$ws.Range("H26").Value = "Tropicales y subtropicales"
$ws.Range("I26").Value = 100108004
$ws.Range("J26").Value = "Papaya"
$ws.Range("K26").Value = "Cultivar IV Región"
$ws.Range("L26").Value = "Primera"
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 17000
$ws.Range("O26").Value = 17000
$ws.Range("P26").Value = 17000
$ws.Range("Q26").Value = "$/bandeja 10 kilos"
$ws.Range("R26").Value = "Provincia del Elquí"
$ws.Range("S26").Value = 1700
$ws.Range("T26").Value = 10

# Row 27
$ws.Range("A27").Value = 3
$ws.Range("B27").Value = "Femacal de La Calera"
$ws.Range("C27").Value = "Coquimbo"
$ws.Range("D27").Value = 44398
$ws.Range("E27").Value = 5
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100108
$ws.Range("H27").Value = "Tropicales y subtropicales"
$ws.Range("I27").Value = 100108004
$ws.Range("J27").Value = "Papaya"
$ws.Range("K27").Value = "Cultivar IV Región"
$ws.Range("L27").Value = "Segunda"
$ws.Range("M27").Value = 50
$ws.Range("N27").Value = 15000
$ws.Range("O27").Value = 15000
$ws.Range("P27").Value = 15000
$ws.Range("Q27").Value = "$/bandeja 10 kilos"
$ws.Range("R27").Value = "Provincia del Elquí"
$ws.Range("S27").Value = 1500
$ws.Range("T27").Value = 10

# Apply the same date display format as the other date cells in column D
$ws.Range("D24:D27").NumberFormat = $ws.Range("D19").NumberFormat
